$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1619
$ws1.Range("F5").Value = 9298
$ws1.Range("F9").Value = 680
$ws1.Range("F17").Value = 1347
$ws1.Range("F20").Value = 1428
$ws1.Range("F21").Value = 105
$ws1.Range("F22").Value = 267
$ws1.Range("F25").Value = 81
$ws1.Range("F40").Value = 169
$ws1.Range("F41").Value = 141

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 3
$ws2.Range("F25").Value = 269
$ws2.Range("F30").Value = 117

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F7").Value = 2274

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1619
$ws4.Range("F5").Value = 9298
$ws4.Range("F10").Value = 680
$ws4.Range("F11").Value = 3
$ws4.Range("F17").Value = 1347
$ws4.Range("F20").Value = 1428
$ws4.Range("F21").Value = 105
$ws4.Range("F22").Value = 267
$ws4.Range("F25").Value = 81
$ws4.Range("F43").Value = 269
$ws4.Range("F44").Value = 169
$ws4.Range("F45").Value = 141
